# edit.ps1 - apply godisnji.docx changes via Word COM-interop
#
# Helper: force a run-split boundary between [from,to) and its neighbours
# by toggling a character-formatting property on/off on that sub-range.
# This does not alter the actual visual formatting (since the value is
# restored), but Word's engine will not re-merge the run back with runs
# that had an intervening explicit formatting operation.
function Mark($doc, $from, $to) {
    $r = $doc.Range($from, $to)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) KLASA:113-02/24-01/  ->  "KLASA:" + " " + "113-02/24-01/"
# -----------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("KLASA:113-02/24-01/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End

$r1 = $d.Range($start, $start + 6)
$r1.InsertAfter(" ")

$p1 = $start + 6
$p2 = $start + 7
$p3 = $end + 1
Mark $d $start $p1
Mark $d $p1 $p2
Mark $d $p2 $p3

# -----------------------------------------------------------------
# 2) URBROJ:2137-37-24-1  ->  "URBROJ:" + " " + "2137-37-24-1"
# -----------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("URBROJ:2137-37-24-1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End

$r1 = $d.Range($start, $start + 7)
$r1.InsertAfter(" ")

$p1 = $start + 7
$p2 = $start + 8
$p3 = $end + 1
Mark $d $start $p1
Mark $d $p1 $p2
Mark $d $p2 $p3

# -----------------------------------------------------------------
# 3) Remove stray spaces before commas in the "Na temelju" sentence
# -----------------------------------------------------------------
$d.Content.Find.Execute("93/14. , 127/17. , 98/19.", $true, $false, $false, $false, $false, $true, 1, $false, "93/14., 127/17., 98/19.", 2) | Out-Null

# -----------------------------------------------------------------
# 4) Remove comma after "Đurđevac" right before "donosi:"
# -----------------------------------------------------------------
$d.Content.Find.Execute("Đurđevac, donosi:", $true, $false, $false, $false, $false, $true, 1, $false, "Đurđevac donosi:", 2) | Out-Null

# -----------------------------------------------------------------
# 5) "Ivi Beljan pripada ..." -> "{{ " + "im_pr_D" + " }} " + "pripada ..."
# -----------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("Ivi Beljan pripada pravo na godišnji odmor za 2024. godinu u trajanju od 30 radnih dana. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End

$r1 = $d.Range($start, $start + 11)
$r1.Text = "{{ im_pr_D }} "

$q0 = $start
$q1 = $start + 3
$q2 = $start + 10
$q3 = $start + 14
$q4 = $end + 3
Mark $d $q0 $q1
Mark $d $q1 $q2
Mark $d $q2 $q3
Mark $d $q3 $q4

# -----------------------------------------------------------------
# 6) Add line spacing (1.5 lines / 360 twips auto) to the "Dostaviti:" paragraph
# -----------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "Dostaviti:" + [char]13) {
        $p.Format.LineSpacingRule = 5
        $p.Format.LineSpacing = 18
    }
}

# -----------------------------------------------------------------
# 7) "2. Tajništvo" -> "2. Tajništv" + "u"
# -----------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("2. Tajništvo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End
$lastCharStart = $end - 1
$rLast = $d.Range($lastCharStart, $end)
$rLast.Text = "u"
Mark $d $start $lastCharStart
Mark $d $lastCharStart $end

# -----------------------------------------------------------------
# 8) "3. Računovodstvo" -> "3. Računovodstv" + "u"
# -----------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("3. Računovodstvo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End
$lastCharStart = $end - 1
$rLast = $d.Range($lastCharStart, $end)
$rLast.Text = "u"
Mark $d $start $lastCharStart
Mark $d $lastCharStart $end
